$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-15 (A..I), newest orders first
$data = @(
    @(491, 45733.22928240741, "Ajay Francis Anchan", 19, 120, 1.8, 1.8, 3.6, "Butterscotch Lassi (x2)"),
    @(490, 45733.22928240741, "Ajay Francis Anchan", 18, 1340, 16.73, 16.73, 35.25, "Chicken Burger (x1), Chicken Cheese Burger (x1), Vanilla Shake (x1), Oreo Shake (x1), Mango Lassi (x1), Strawberry Lassi (x1), Butterscotch Lassi (x1), Chicken Wrap (x1), Veg Wrap (x4), Chicken Cheese Pops (x1), Veg Cheese Pops (x3), Banana Shake (x2), Black Currant Shake (x1), Strawberry Shake (x1), Chocolate Shake (x1), Hard Rock Coffee (x1), 8PM Coffee (x1), Coffee Italia (x1), Belgian Coffee (x1), Lime Juice (x1)"),
    @(489, 45733.22928240741, "Ajay Francis Anchan", 17, 330, 8, 8, 16, "Strawberry Lassi (x1), Chicken Wrap (x4), Banana Shake (x1)"),
    @(488, 45733.22928240741, "Ajay Francis Anchan", 16, 190, 5, 5, 10, "Strawberry Lassi (x1), Chocolate Shake (x1), Hard Rock Coffee (x1), Belgian Coffee (x1), Lime Juice (x1)"),
    @(487, 45733.22928240741, "Ajay Francis Anchan", 15, 100, 1.25, 1.25, 2.5, "Oreo Shake (x2)"),
    @(486, 45733.22928240741, "Ajay Francis Anchan", 14, 210, 1.93, 1.93, 5.65, "Chicken Burger (x1), Chicken Cheese Burger (x1), Vanilla Shake (x1), Oreo Shake (x1)"),
    @(485, 45733.22928240741, "Ajay Francis Anchan", 13, 170, 2.42, 2.42, 4.85, "Mango Lassi (x2), Veg Wrap (x1)"),
    @(484, 45733.22928240741, "Ajay Francis Anchan", 12, 110, 1.4, 1.4, 2.8, "Butterscotch Lassi (x1), Banana Shake (x1)"),
    @(483, 45733.22928240741, "Ajay Francis Anchan", 11, 140, 2.45, 2.45, 4.9, "Chicken Cheese Pops (x1), Veg Cheese Pops (x1)"),
    @(482, 45733.22928240741, "Ajay Francis Anchan", 10, 140, 1.7, 1.7, 3.4, "Strawberry Lassi (x2), Butterscotch Lassi (x1)"),
    @(481, 45733.22928240741, "Ajay Francis Anchan", 9, 200, 4, 4, 8, "Banana Shake (x1), Black Currant Shake (x2), Strawberry Shake (x1)"),
    @(480, 45733.22928240741, "Ajay Francis Anchan", 8, 300, 5.75, 5.75, 11.5, "Chicken Cheese Burger (x1), Chicken Wrap (x2), Veg Cheese Pops (x1), Banana Shake (x1)"),
    @(478, 45733.22928240741, "Ajay Francis Anchan", 6, 70, 1.23, 1.23, 2.45, "Chicken Cheese Pops (x1)"),
    @(477, 45733.22928240741, "Ajay Francis Anchan", 5, 120, 1.73, 1.73, 3.45, "Veg Cheese Pops (x1), Banana Shake (x1)")
)

# New rows (7-15) need the date style (from B2) copied onto column B
# before the value is written, so the number format matches the
# existing date column instead of defaulting to General.
$ws.Range("B2").Copy()
for ($r = 7; $r -le 15; $r++) {
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)  # xlPasteFormats
}

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}